# Update latest output (run 163)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) columns ---
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("E2").Value = 369.6368129999999
$wsSchedule.Range("F2").Value = 8.148959722222219
$wsSchedule.Range("E3").Value = 416.68462875
$wsSchedule.Range("F3").Value = 27.55850719246032

# --- Sheet "Detailed": update Price column (B) and some Type column (C) values ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B13").Value = 78
$wsDetailed.Range("B14").Value = 64.89

$wsDetailed.Range("B15").Value = 59.96022
$wsDetailed.Range("C15").Value = "historical"

$wsDetailed.Range("B16").Value = 51.82055
$wsDetailed.Range("C16").Value = "historical"

$wsDetailed.Range("B17").Value = 0.66384
$wsDetailed.Range("B18").Value = -5.05758
$wsDetailed.Range("B19").Value = -6.11025
$wsDetailed.Range("B20").Value = -6.79343
$wsDetailed.Range("B21").Value = -7.70579
$wsDetailed.Range("B22").Value = -7.53666
$wsDetailed.Range("B23").Value = -8.276339999999999
$wsDetailed.Range("B24").Value = -7.58016
$wsDetailed.Range("B25").Value = -5.58973
$wsDetailed.Range("B26").Value = -6.49292
$wsDetailed.Range("B27").Value = -6.65905
$wsDetailed.Range("B28").Value = -7.19834
$wsDetailed.Range("B29").Value = -6.07345
$wsDetailed.Range("B31").Value = -0.88236
$wsDetailed.Range("B32").Value = -2.54301
$wsDetailed.Range("B33").Value = 0.00001
$wsDetailed.Range("B34").Value = -9.709580000000001
$wsDetailed.Range("B37").Value = -7.25528
$wsDetailed.Range("B38").Value = -0.39875
$wsDetailed.Range("B39").Value = 7.21234
$wsDetailed.Range("B40").Value = 18.95614
$wsDetailed.Range("B42").Value = 55.33037
$wsDetailed.Range("B43").Value = 51.4753
$wsDetailed.Range("B44").Value = 57.01318
$wsDetailed.Range("B46").Value = 43.61794
$wsDetailed.Range("B49").Value = 48.93665
